$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns (and the B/C cells we touch) to Text format
# before writing, so numeric-looking strings like "0.636" or "58.07" are not
# auto-converted to numbers by Excel and keep their exact textual form
# (trailing zeros, multi-dot thousand separators, etc.), matching the
# original inline-string cells.
$ws.Range("B43:E44").NumberFormat = "@"
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "37.601.24"
$ws.Range("E2").Value = "  -0.29%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.085.07"
$ws.Range("E3").Value = "  +0.45%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "233.72"

# Row 6 - XRP
$ws.Range("D6").Value = "0.636"
$ws.Range("E6").Value = "  +2.19%  "

# Row 8 - Solana
$ws.Range("D8").Value = "58.07"
$ws.Range("E8").Value = "  -0.14%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.392"
$ws.Range("E9").Value = "  -1.01%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0780"
$ws.Range("E10").Value = "  -0.41%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +3.02%  "

# Row 12 - Chainlink
$ws.Range("D12").Value = "15.12"
$ws.Range("E12").Value = "  +2.41%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.394.83"
$ws.Range("E13").Value = "  +0.51%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "21.16"
$ws.Range("E14").Value = "  +1.50%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.27%  "

# Row 16 - Polkadot
$ws.Range("D16").Value = "5.36"
$ws.Range("E16").Value = "  +1.03%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.084.21"
$ws.Range("E17").Value = "  -0.14%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "37.576.88"
$ws.Range("E18").Value = "  -0.13%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  -1.88%  "

# Row 20 - Litecoin
$ws.Range("D20").Value = "70.84"
$ws.Range("E20").Value = "  -0.30%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0834"
$ws.Range("E21").Value = "  +0.17%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "229.26"
$ws.Range("E22").Value = "  +0.39%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.09%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "2.37"
$ws.Range("E24").Value = "  -1.03%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -0.31%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "9.70"
$ws.Range("E26").Value = "  +7.29%  "

# Row 27 - Monero
$ws.Range("D27").Value = "170.72"
$ws.Range("E27").Value = "  +0.05%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  -3.79%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "19.57"
$ws.Range("E29").Value = "  +0.67%  "

# Row 30 - ImmutableX
$ws.Range("D30").Value = "1.38"
$ws.Range("E30").Value = "  +0.65%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "0.122"
$ws.Range("E31").Value = "  +1.02%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "4.66"
$ws.Range("E32").Value = "  +0.03%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.0638"

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").Value = "4.66"
$ws.Range("E34").Value = "  +0.73%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").Value = "2.49"
$ws.Range("E35").Value = "  +0.23%  "

# Row 36 - WEMIXToken
$ws.Range("E36").Value = "  -0.81%  "

# Row 37 - RenderToken
$ws.Range("E37").Value = "  -1.63%  "

# Row 38 - BinanceUSD
$ws.Range("E38").Value = "  -0.06%  "

# Row 39 - THORChain
$ws.Range("D39").Value = "5.38"

# Row 40 - VeChain
$ws.Range("D40").Value = "0.0232"
$ws.Range("E40").Value = "  +8.47%  "

# Row 41 - Aave
$ws.Range("D41").Value = "100.79"
$ws.Range("E41").Value = "  +2.84%  "

# Row 42 - Cronos
$ws.Range("E42").Value = "  -0.78%  "

# Row 43 - now HuobiToken (was TrustWalletToken)
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "2.90"
$ws.Range("E43").Value = "  +0.66%  "

# Row 44 - now TrustWalletToken (was HuobiToken)
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "1.19"
$ws.Range("E44").Value = "  +3.63%  "

# Row 45 - InjectiveProtocol
$ws.Range("D45").Value = "16.87"
$ws.Range("E45").Value = "  +3.10%  "

# Row 46 - Maker
$ws.Range("D46").Value = "1.456.93"
$ws.Range("E46").Value = "  +0.66%  "

# Row 47 - ARBITRUM
$ws.Range("E47").Value = "  -0.58%  "

# Row 48 - FTXToken
$ws.Range("D48").Value = "4.00"
$ws.Range("E48").Value = "  -5.47%  "

# Row 49 - FraxShare
$ws.Range("D49").Value = "7.25"
$ws.Range("E49").Value = "  -2.21%  "

# Row 50 - MXToken
$ws.Range("D50").Value = "2.95"

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.278.55"
$ws.Range("E51").Value = "  +0.49%  "

# Restore the default ("Normal") style on the touched ranges so the cells'
# style index matches the original workbook (no explicit style attribute),
# now that the text has been safely written as Text-formatted values.
$ws.Range("B43:E44").Style = "Normal"
$ws.Range("D2:E51").Style = "Normal"
